{"js": "// The document contains a table of \"three-digit \u00f7 one-digit\" division\n// problems whose worked-out answers were regenerated. Each OLD answer\n// string is unique in the document, so body.search(old, {matchCase:true})\n// unambiguously locates the single run to rewrite with the NEW answer.\nconst replacements = [\n  [\"637\u00f75=127, 2\", \"343\u00f77=49, 0\"],\n  [\"735\u00f72=367, 1\", \"685\u00f76=114, 1\"],\n  [\"868\u00f72=434, 0\", \"985\u00f75=197, 0\"],\n  [\"577\u00f77=82, 3\", \"192\u00f72=96, 0\"],\n  [\"885\u00f73=295, 0\", \"995\u00f72=497, 1\"],\n  [\"638\u00f78=79, 6\", \"575\u00f72=287, 1\"],\n  [\"599\u00f79=66, 5\", \"404\u00f78=50, 4\"],\n  [\"796\u00f74=199, 0\", \"462\u00f78=57, 6\"],\n  [\"666\u00f73=222, 0\", \"997\u00f78=124, 5\"],\n  [\"566\u00f74=141, 2\", \"459\u00f75=91, 4\"],\n  [\"303\u00f72=151, 1\", \"658\u00f77=94, 0\"],\n  [\"879\u00f75=175, 4\", \"644\u00f78=80, 4\"],\n  [\"470\u00f77=67, 1\", \"545\u00f76=90, 5\"],\n  [\"159\u00f77=22, 5\", \"948\u00f72=474, 0\"],\n  [\"414\u00f77=59, 1\", \"172\u00f75=34, 2\"],\n  [\"144\u00f74=36, 0\", \"920\u00f77=131, 3\"],\n  [\"507\u00f77=72, 3\", \"391\u00f76=65, 1\"],\n  [\"248\u00f72=124, 0\", \"315\u00f74=78, 3\"],\n  [\"547\u00f77=78, 1\", \"533\u00f76=88, 5\"],\n  [\"370\u00f77=52, 6\", \"524\u00f72=262, 0\"],\n  [\"167\u00f78=20, 7\", \"925\u00f79=102, 7\"],\n  [\"516\u00f72=258, 0\", \"747\u00f72=373, 1\"],\n  [\"879\u00f76=146, 3\", \"790\u00f75=158, 0\"],\n  [\"877\u00f79=97, 4\", \"292\u00f72=146, 0\"],\n  [\"322\u00f79=35, 7\", \"101\u00f74=25, 1\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the old three-digit-division answers with the new ones.\n# Each old value is unique within the document, so Find/Replace against the\n# whole document body (Content) unambiguously targets the single matching run.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"637\u00f75=127, 2\", \"343\u00f77=49, 0\"),\n    @(\"735\u00f72=367, 1\", \"685\u00f76=114, 1\"),\n    @(\"868\u00f72=434, 0\", \"985\u00f75=197, 0\"),\n    @(\"577\u00f77=82, 3\", \"192\u00f72=96, 0\"),\n    @(\"885\u00f73=295, 0\", \"995\u00f72=497, 1\"),\n    @(\"638\u00f78=79, 6\", \"575\u00f72=287, 1\"),\n    @(\"599\u00f79=66, 5\", \"404\u00f78=50, 4\"),\n    @(\"796\u00f74=199, 0\", \"462\u00f78=57, 6\"),\n    @(\"666\u00f73=222, 0\", \"997\u00f78=124, 5\"),\n    @(\"566\u00f74=141, 2\", \"459\u00f75=91, 4\"),\n    @(\"303\u00f72=151, 1\", \"658\u00f77=94, 0\"),\n    @(\"879\u00f75=175, 4\", \"644\u00f78=80, 4\"),\n    @(\"470\u00f77=67, 1\", \"545\u00f76=90, 5\"),\n    @(\"159\u00f77=22, 5\", \"948\u00f72=474, 0\"),\n    @(\"414\u00f77=59, 1\", \"172\u00f75=34, 2\"),\n    @(\"144\u00f74=36, 0\", \"920\u00f77=131, 3\"),\n    @(\"507\u00f77=72, 3\", \"391\u00f76=65, 1\"),\n    @(\"248\u00f72=124, 0\", \"315\u00f74=78, 3\"),\n    @(\"547\u00f77=78, 1\", \"533\u00f76=88, 5\"),\n    @(\"370\u00f77=52, 6\", \"524\u00f72=262, 0\"),\n    @(\"167\u00f78=20, 7\", \"925\u00f79=102, 7\"),\n    @(\"516\u00f72=258, 0\", \"747\u00f72=373, 1\"),\n    @(\"879\u00f76=146, 3\", \"790\u00f75=158, 0\"),\n    @(\"877\u00f79=97, 4\", \"292\u00f72=146, 0\"),\n    @(\"322\u00f79=35, 7\", \"101\u00f74=25, 1\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $found = $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n\n    if (-not $found) {\n        throw \"Text not found: $old\"\n    }\n}\n"}
